# Update sentencing dates in the Judgment Entry document.
$d = $word.ActiveDocument

# Replace every occurrence of "June 12, 2022" with "June 13, 2022".
# This covers:
#   - " on June 12, 2022."
#   - "June 12, 2022" (bold date line)
#   - " license is suspended from June 12, 2022"
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("June 12, 2022", $true, $false, $false, $false, $false, `
               $true, 1, $false, "June 13, 2022", 2)

# Replace "August 11, 2022" with "August 12, 2022".
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute("August 11, 2022", $true, $false, $false, $false, $false, `
                $true, 1, $false, "August 12, 2022", 2)
